# Trade #26 (MarketMaking trade #56, row 57 on "All Trades" / row 28 on
# "MarketMaking") closes, and a brand-new OPEN "momentum" trade (#85) is
# logged at 2026-02-18 00:13:02. Summary / Strategy Status roll-ups are
# updated to match the new totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while keeping it as literal TEXT
# (prevents the engine's auto date/number detection from turning things
# like "2026-02-18" into a date serial), then restore the cell style so
# no stray number-format residue is left behind.
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.78   # Current Capital
$summary.Range("B4").Value = 0.88      # Total P&L $
$summary.Range("B5").Value = 0.33      # Total P&L %
$summary.Range("B6").Value = 54        # Total Trades
$summary.Range("B8").Value = 21        # Losing Trades
$summary.Range("B9").Value = 55.56     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet, MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.78
$status.Range("D6").Value = 25
$status.Range("E6").Value = -0.03
$status.Range("F6").Value = -0.22
$status.Range("G6").Value = 56

# ---------------------------------------------------------------------
# 3) "All Trades" sheet - close trade row 57 (Trade #56, MarketMaking)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G57").Value = 0.69697
$allTrades.Range("H57").Value = "CLOSED"
$allTrades.Range("I57").Value = -10.6449
$allTrades.Range("J57").Value = -0.08
$allTrades.Range("K57").Value = 99.78
$allTrades.Range("L57").Value = "early_exit"
$allTrades.Range("M57").Value = 0.12

# Append the new OPEN momentum trade as row 86
Set-TextValue $allTrades.Range("B86") "2026-02-18"
Set-TextValue $allTrades.Range("C86") "00:13:02"
$allTrades.Range("A86").Value = 85
Set-TextValue $allTrades.Range("D86") "momentum"
Set-TextValue $allTrades.Range("E86") "DOWN"
$allTrades.Range("F86").Value = 0.78
Set-TextValue $allTrades.Range("H86") "OPEN"
$allTrades.Range("I86").Value = 0
$allTrades.Range("J86").Value = 0
$allTrades.Range("K86").Value = 100
$allTrades.Range("M86").Value = 0
$allTrades.Range("N86").Value = 0
$allTrades.Range("O86").Value = 0
$allTrades.Range("P86").Value = 0.9
Set-TextValue $allTrades.Range("Q86") "Downward momentum: -1.980% over 10 samples"

# ---------------------------------------------------------------------
# 4) "momentum" strategy sheet - append the same new trade as row 16
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A16").Value = 85
Set-TextValue $momentum.Range("B16") "2026-02-18"
Set-TextValue $momentum.Range("C16") "00:13:02"
Set-TextValue $momentum.Range("D16") "momentum"
Set-TextValue $momentum.Range("E16") "DOWN"
$momentum.Range("F16").Value = 0.78
Set-TextValue $momentum.Range("H16") "OPEN"
$momentum.Range("I16").Value = 0
$momentum.Range("J16").Value = 0
$momentum.Range("K16").Value = 100
$momentum.Range("L16").Value = 0
$momentum.Range("M16").Value = 0
$momentum.Range("N16").Value = 0.9
Set-TextValue $momentum.Range("O16") "Downward momentum: -1.980% over 10 samples"
$momentum.Range("Q16").Value = 0

# ---------------------------------------------------------------------
# 5) "MarketMaking" strategy sheet - close the same trade as row 28
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G28").Value = 0.69697
$marketMaking.Range("H28").Value = "CLOSED"
$marketMaking.Range("I28").Value = -10.6449
$marketMaking.Range("J28").Value = -0.08
$marketMaking.Range("K28").Value = 99.78
$marketMaking.Range("P28").Value = "early_exit"
$marketMaking.Range("Q28").Value = 0.12
